## Adding support for multiple users on a computer
## phy-tglab11 (column G) is now shared by two users (Hernan/hgarcia and
## Eric/emsthree), each with their own set of folders. We insert a new
## column for the second user's folders, and a new header row that
## records which user name each computer's column belongs to.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column after the existing "phy-tglab11" column (G) to hold
# Eric's (emsthree) folder paths for that same machine, and a new row
# right under the header to hold the "User Name" for each computer.
$ws.Columns("H").Insert()
$ws.Rows(2).Insert()

# Header row: the new column H is the same computer (phy-tglab11) as G.
$ws.Range("H1").Value = "phy-tglab11"

# New "User Name" row distinguishing the two users of phy-tglab11.
$ws.Range("A2").Value = "User Name"
$ws.Range("G2").Value = "hgarcia"
$ws.Range("H2").Value = "emsthree"

# Eric's (emsthree) folders on phy-tglab11.
$ws.Range("H3").Value = "F:\Eric\Local\[07] Transcription\RawData"
$ws.Range("H4").Value = "F:\Eric\Local\[07] Transcription\FISHAnalysisData"
$ws.Range("H5").Value = "F:\Eric\Dropbox\LivemRNAData"
$ws.Range("H8").Value = "F:\Eric\GitHub\mRNADynamics"

# Match the tweaked column widths from the authored workbook.
$ws.Columns("H").ColumnWidth = 44.5

# Restore the split view over the new layout (pane boundary now sits
# after column G, i.e. before the newly inserted column H) and leave the
# selection on the newly added cells.
$ws.Activate()
$excel.ActiveWindow.SplitColumn = 6
$excel.ActiveWindow.SplitRow = 0
$ws.Range("H2").Select()
